# Auto-generated edit script: updates Universalis market price / leve profit
# figures across multiple job sheets (ALC, ARM, BSM, CRP, LTW, WVR) to match
# the latest scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 62
$ws_ALC.Range("H62").Value = 16673667
$ws_ALC.Range("I62").Value = 16673667
$ws_ALC.Range("K62").Value = 16673667
$ws_ALC.Range("M62").Value = -16673043

# ALC row 65
$ws_ALC.Range("H65").Value = 16673667
$ws_ALC.Range("I65").Value = 16673667
$ws_ALC.Range("K65").Value = 83368335
$ws_ALC.Range("M65").Value = -83365215

# ALC row 112
$ws_ALC.Range("H112").Value = 84924.75
$ws_ALC.Range("J112").Value = 145140.42
$ws_ALC.Range("L112").Value = 435421.26
$ws_ALC.Range("N112").Value = -437637.26

# ARM row 32
$ws_ARM.Range("H32").Value = 3777.2388
$ws_ARM.Range("I32").Value = 2858
$ws_ARM.Range("K32").Value = 2858
$ws_ARM.Range("M32").Value = -2571

# ARM row 74
$ws_ARM.Range("H74").Value = 5655.7
$ws_ARM.Range("I74").Value = 1143.0294
$ws_ARM.Range("K74").Value = 1143.0294
$ws_ARM.Range("M74").Value = -269.0293999999999

# ARM row 77
$ws_ARM.Range("H77").Value = 5655.7
$ws_ARM.Range("I77").Value = 1143.0294
$ws_ARM.Range("K77").Value = 5715.146999999999
$ws_ARM.Range("M77").Value = -1347.146999999999

# ARM row 132
$ws_ARM.Range("H132").Value = 2515.8462
$ws_ARM.Range("I132").Value = 1822
$ws_ARM.Range("K132").Value = 5466
$ws_ARM.Range("M132").Value = -2936

# BSM row 64
$ws_BSM.Range("H64").Value = 947
$ws_BSM.Range("J64").Value = 701.5
$ws_BSM.Range("L64").Value = 701.5
$ws_BSM.Range("N64").Value = -1151.5

# BSM row 67
$ws_BSM.Range("H67").Value = 947
$ws_BSM.Range("J67").Value = 701.5
$ws_BSM.Range("L67").Value = 701.5
$ws_BSM.Range("N67").Value = -2261.5

# BSM row 86
$ws_BSM.Range("H86").Value = 4457.697
$ws_BSM.Range("I86").Value = 1552.8518
$ws_BSM.Range("J86").Value = 17529.5
$ws_BSM.Range("K86").Value = 1552.8518
$ws_BSM.Range("L86").Value = 17529.5
$ws_BSM.Range("M86").Value = -429.8517999999999
$ws_BSM.Range("N86").Value = -19775.5

# BSM row 89
$ws_BSM.Range("H89").Value = 4457.697
$ws_BSM.Range("I89").Value = 1552.8518
$ws_BSM.Range("J89").Value = 17529.5
$ws_BSM.Range("K89").Value = 7764.259
$ws_BSM.Range("L89").Value = 87647.5
$ws_BSM.Range("M89").Value = -2148.259
$ws_BSM.Range("N89").Value = -98879.5

# BSM row 134
$ws_BSM.Range("H134").Value = 1948.7222
$ws_BSM.Range("I134").Value = 1960.4117
$ws_BSM.Range("K134").Value = 5881.2351
$ws_BSM.Range("M134").Value = -3346.2351

# CRP row 31
$ws_CRP.Range("H31").Value = 47323.22
$ws_CRP.Range("I31").Value = 51448.25
$ws_CRP.Range("J31").Value = 19823
$ws_CRP.Range("K31").Value = 51448.25
$ws_CRP.Range("L31").Value = 19823
$ws_CRP.Range("M31").Value = -51153.25
$ws_CRP.Range("N31").Value = -20413

# CRP row 34
$ws_CRP.Range("H34").Value = 47323.22
$ws_CRP.Range("I34").Value = 51448.25
$ws_CRP.Range("J34").Value = 19823
$ws_CRP.Range("K34").Value = 51448.25
$ws_CRP.Range("L34").Value = 19823
$ws_CRP.Range("M34").Value = -51246.25
$ws_CRP.Range("N34").Value = -20227

# CRP row 58
$ws_CRP.Range("H58").Value = 2297.4614
$ws_CRP.Range("I58").Value = 2388.75
$ws_CRP.Range("K58").Value = 2388.75
$ws_CRP.Range("M58").Value = -2185.75

# CRP row 99
$ws_CRP.Range("H99").Value = 8598.134
$ws_CRP.Range("I99").Value = 8181.3335
$ws_CRP.Range("J99").Value = 9223.333000000001
$ws_CRP.Range("K99").Value = 8181.3335
$ws_CRP.Range("L99").Value = 9223.333000000001
$ws_CRP.Range("M99").Value = -6683.3335
$ws_CRP.Range("N99").Value = -12219.333

# CRP row 126
$ws_CRP.Range("H126").Value = 8598.134
$ws_CRP.Range("I126").Value = 8181.3335
$ws_CRP.Range("J126").Value = 9223.333000000001
$ws_CRP.Range("K126").Value = 24544.0005
$ws_CRP.Range("L126").Value = 27669.999
$ws_CRP.Range("M126").Value = -22074.0005
$ws_CRP.Range("N126").Value = -32609.999

# CRP row 132
$ws_CRP.Range("H132").Value = 4353.5356
$ws_CRP.Range("I132").Value = 4280.731
$ws_CRP.Range("K132").Value = 12842.193
$ws_CRP.Range("M132").Value = -10312.193

# CRP row 136
$ws_CRP.Range("H136").Value = 2297.4614
$ws_CRP.Range("I136").Value = 2388.75
$ws_CRP.Range("K136").Value = 7166.25
$ws_CRP.Range("M136").Value = -4616.25

# CRP row 138
$ws_CRP.Range("H138").Value = 55390
$ws_CRP.Range("I138").Value = 30000
$ws_CRP.Range("K138").Value = 30000
$ws_CRP.Range("M138").Value = -24860

# LTW row 16
$ws_LTW.Range("H16").Value = 11999.75
$ws_LTW.Range("I16").Value = 0
$ws_LTW.Range("J16").Value = 11999.75
$ws_LTW.Range("K16").Value = 0
$ws_LTW.Range("L16").Value = 11999.75
$ws_LTW.Range("M16").ClearContents()
$ws_LTW.Range("N16").Value = -12339.75

# LTW row 22
$ws_LTW.Range("H22").Value = 1563.2142
$ws_LTW.Range("I22").Value = 1837.5
$ws_LTW.Range("J22").Value = 1517.5
$ws_LTW.Range("K22").Value = 1837.5
$ws_LTW.Range("L22").Value = 1517.5
$ws_LTW.Range("M22").Value = -1542.5
$ws_LTW.Range("N22").Value = -2107.5

# LTW row 27
$ws_LTW.Range("H27").Value = 1563.2142
$ws_LTW.Range("I27").Value = 1837.5
$ws_LTW.Range("J27").Value = 1517.5
$ws_LTW.Range("K27").Value = 1837.5
$ws_LTW.Range("L27").Value = 1517.5
$ws_LTW.Range("M27").Value = -1730.5
$ws_LTW.Range("N27").Value = -1731.5

# LTW row 46
$ws_LTW.Range("H46").Value = 1681.8462
$ws_LTW.Range("I46").Value = 1538.9
$ws_LTW.Range("K46").Value = 1538.9
$ws_LTW.Range("M46").Value = -1350.9

# LTW row 61
$ws_LTW.Range("H61").Value = 3978.125
$ws_LTW.Range("I61").Value = 3978.125
$ws_LTW.Range("K61").Value = 3978.125
$ws_LTW.Range("M61").Value = -3776.125

# LTW row 113
$ws_LTW.Range("H113").Value = 3978.125
$ws_LTW.Range("I113").Value = 3978.125
$ws_LTW.Range("K113").Value = 3978.125
$ws_LTW.Range("M113").Value = -1808.125

# LTW row 132
$ws_LTW.Range("H132").Value = 4270.2
$ws_LTW.Range("I132").Value = 3700.4443
$ws_LTW.Range("J132").Value = 5124.8335
$ws_LTW.Range("K132").Value = 11101.3329
$ws_LTW.Range("L132").Value = 15374.5005
$ws_LTW.Range("M132").Value = -8571.332900000001
$ws_LTW.Range("N132").Value = -20434.5005

# LTW row 136
$ws_LTW.Range("H136").Value = 3768.6191
$ws_LTW.Range("I136").Value = 3472.6
$ws_LTW.Range("J136").Value = 4508.6665
$ws_LTW.Range("K136").Value = 10417.8
$ws_LTW.Range("L136").Value = 13525.9995
$ws_LTW.Range("M136").Value = -7867.799999999999
$ws_LTW.Range("N136").Value = -18625.9995

# WVR row 123 (quantity dropped to 0, so NQ/HQ leve profit no longer applies;
# clear the stale HQ-profit figure (N123) to match the now-zeroed price/qty
# inputs for this row, mirroring the all-zero pattern used elsewhere in this
# sheet, e.g. rows 120-121)
$ws_WVR.Range("H123").Value = 0
$ws_WVR.Range("J123").Value = 0
$ws_WVR.Range("L123").Value = 0
$ws_WVR.Range("M123").ClearContents()
$ws_WVR.Range("N123").ClearContents()

# WVR row 126
$ws_WVR.Range("H126").Value = 3241.75
$ws_WVR.Range("I126").Value = 2590.0625
$ws_WVR.Range("J126").Value = 5848.5
$ws_WVR.Range("K126").Value = 7770.1875
$ws_WVR.Range("L126").Value = 17545.5
$ws_WVR.Range("M126").Value = -5300.1875
$ws_WVR.Range("N126").Value = -22485.5

# WVR row 132
$ws_WVR.Range("H132").Value = 1973.5
$ws_WVR.Range("I132").Value = 2105.4285
$ws_WVR.Range("J132").Value = 1050
$ws_WVR.Range("K132").Value = 6316.2855
$ws_WVR.Range("L132").Value = 3150
$ws_WVR.Range("M132").Value = -3786.2855
$ws_WVR.Range("N132").Value = -8210

# WVR row 136
$ws_WVR.Range("H136").Value = 2382.5
$ws_WVR.Range("I136").Value = 2189.842
$ws_WVR.Range("K136").Value = 6569.526
$ws_WVR.Range("M136").Value = -4019.526
